# Update "paises.xlsx" (sheet "Pais") with refreshed COVID-19 country data
# and swap two pairs of countries whose case counts re-ordered them in the
# (descending, by total cases) table:
#   - Row 37 / Row 38: Colombia / Polonia swap places (Polonia's new total
#     of 22074 overtakes Colombia's 21981).
#   - Row 207 / Row 208: Groenlandia / Islas Turcas y Caicos swap places.
# Also refreshes the "Datos actualizados..." timestamp in cell A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title cell: refresh the "last updated" timestamp (18:05 -> 18:35)
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 18:35"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1711569
$ws.Cells.Item(4, 3).Value = 5343
$ws.Cells.Item(4, 4).Value = 467962
$ws.Cells.Item(4, 5).Value = 1143620
$ws.Cells.Item(4, 7).Value = 182
$ws.Cells.Item(4, 8).Value = 99987

# Row 37: was Colombia, now Polonia (moved up: 22074 total cases)
$ws.Cells.Item(37, 1).Value = "Polonia"
$ws.Cells.Item(37, 2).Value = 22074
$ws.Cells.Item(37, 3).Value = 443
$ws.Cells.Item(37, 4).Value = 10020
$ws.Cells.Item(37, 5).Value = 11030
$ws.Cells.Item(37, 7).Value = 17
$ws.Cells.Item(37, 8).Value = 1024

# Row 38: was Polonia, now Colombia (moved down: 21981 total cases)
$ws.Cells.Item(38, 1).Value = "Colombia"
$ws.Cells.Item(38, 2).Value = 21981
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(38, 4).Value = 5265
$ws.Cells.Item(38, 5).Value = 15966
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 750

# Row 54: Chequia
$ws.Cells.Item(54, 2).Value = 9034
$ws.Cells.Item(54, 3).Value = 32
$ws.Cells.Item(54, 4).Value = 6263
$ws.Cells.Item(54, 5).Value = 2454

# Row 69: Irak
$ws.Cells.Item(69, 4).Value = 2852
$ws.Cells.Item(69, 5).Value = 1827

# Row 103: Sri Lanka
$ws.Cells.Item(103, 2).Value = 1317
$ws.Cells.Item(103, 3).Value = 135
$ws.Cells.Item(103, 5).Value = 595

# Row 126: Jordania
$ws.Cells.Item(126, 4).Value = 586
$ws.Cells.Item(126, 5).Value = 123

# Row 144: Isla de Man
$ws.Cells.Item(144, 4).Value = 306
$ws.Cells.Item(144, 5).Value = 6

# Row 207: was Groenlandia, now Islas Turcas y Caicos
$ws.Cells.Item(207, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(207, 4).Value = 10
$ws.Cells.Item(207, 8).Value = 1

# Row 208: was Islas Turcas y Caicos, now Groenlandia
$ws.Cells.Item(208, 1).Value = "Groenlandia"
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(208, 8).Value = 0
